$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 65

$ws.Range("B5").Value = 54
$ws.Range("C5").Value = 7

$ws.Range("B6").Value = 35
$ws.Range("C6").Value = 13

$ws.Range("B7").Value = 36
$ws.Range("C7").Value = 7

$ws.Range("B9").Value = 38
$ws.Range("C9").Value = 23

$ws.Range("B12").Value = 33
$ws.Range("C12").Value = 23

$ws.Range("B13").Value = 22
$ws.Range("C13").Value = 26
